$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new experiment row (row 19) following the same pattern as the
# existing rows (2-18): Local Area Density, Potential Radius,
# Local/Global Inhibition, NumActiveColumnsPerInhArea, Result Image Name.
$ws.Range("A19").Value = "Exp 23"
$ws.Range("B19").Value = 0.25
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = "Local"
$ws.Range("E19").Value = -1
$ws.Range("F19").Value = "Exp 23.png"

# Match the centered-alignment style used by the rest of the data rows
# (columns A:E) for the new row.
$ws.Range("A19:E19").HorizontalAlignment = -4108

# Update the active selection, matching the author's final cursor position.
$ws.Range("J17").Select()
